$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E20").Value = "{`n     ""temperature"": 120.5,`n     ""unit"": ""Celsius"",`n     ""time"": ""2023-07-12T16:21:53.389+02:00"",`n     ""externalId"": ""berlin_01"",`n     ""unexpected"": 17.5`n}"

$ws.Range("G20").Value = "A measasurement should be created for the device berlin_01.`nThe fragment ""c8y_Fragment_to_remove"" is not included in the created measurement, as the repair strategy is ""REMOVE_IF_NULL"".`nIn addition the reapar strategy ""CREATE_IF_MISSING"" is used. Thjsi is required to map the node ""unexpected"" to the target fragment ""c8y_Unexpected"". This is created, due to the used reapir strategy."

$ws.Rows.Item(20).RowHeight = 252
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
